$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously-blank cells on row 18 (Alex Schell)
$ws.Range("B18").Value = "L"
$ws.Range("C18").Value = "N"
$ws.Range("D18").Value = "N"

# Match the formatting ("s=1" cell style - Arial 10, no explicit font
# color) already used by the rest of the row/sheet, by copying it from
# the neighboring A18 cell.
$ws.Range("A18").Copy()
$ws.Range("B18:D18").PasteSpecial(-4122)

# Move the view/selection from the bottom of the sheet back up to the
# newly edited row so D18 is the active cell and the sheet is no longer
# scrolled down to row 13.
$ws.Range("D18").Select()
